$wb = $excel.ActiveWorkbook

# New GUID-named handoff file (replaces 4e6f9cc2-cbad-4b98-b049-9b603a1ddb9d)
# and its new content hash (replaces ef90965287d3ccd24169d5d3b3a45c1037bfbb91),
# generated for this report-for-handoff run.
$newId = "30267ee3-b3ef-4d4b-92d3-00d4538bc873"
$newZhHash = "d259d4139d183ed928138f6ff975a64695e163be"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-08-17 14:58:59"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newId.md"
}

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("G2").Value = "$newId.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-17 14:58:54"
foreach ($h in $wsZh.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("G2").Value = "$newId.$newZhHash.de-de.xlf"
foreach ($h in $wsDe.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}
